# Natmi following Dr Hou advice
# Update LR-pairs data: sending-cluster x target-cluster matrix now includes
# "FAPs" as a target cluster (previously missing), bringing the grid from
# 4 sending clusters x 3 target clusters (12 rows) to 4 x 4 (16 rows).
# All numeric columns (E:T) are refreshed with the recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 162.399297
$ws.Cells.Item(2,8).Value = 487.197891
$ws.Cells.Item(2,9).Value = 0.3910371682630009
$ws.Cells.Item(2,10).Value = 0.3910371682630009
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 21.38398133333333
$ws.Cells.Item(2,14).Value = 64.151944
$ws.Cells.Item(2,15).Value = 0.864235283869651
$ws.Cells.Item(2,16).Value = 0.8642352838696511
$ws.Cells.Item(2,17).Value = 3472.743535594456
$ws.Cells.Item(2,18).Value = 31254.6918203501
$ws.Cells.Item(2,19).Value = 0.3379481181173591
$ws.Cells.Item(2,20).Value = 0.3379481181173591

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 162.399297
$ws.Cells.Item(3,8).Value = 487.197891
$ws.Cells.Item(3,9).Value = 0.3910371682630009
$ws.Cells.Item(3,10).Value = 0.3910371682630009
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.05505166666666667
$ws.Cells.Item(3,14).Value = 0.165155
$ws.Cells.Item(3,15).Value = 0.002224917428963528
$ws.Cells.Item(3,16).Value = 0.002224917428963528
$ws.Cells.Item(3,17).Value = 8.940351965345
$ws.Cells.Item(3,18).Value = 80.46316768810499
$ws.Cells.Item(3,19).Value = 0.0008700254110408944
$ws.Cells.Item(3,20).Value = 0.0008700254110408946

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 162.399297
$ws.Cells.Item(4,8).Value = 487.197891
$ws.Cells.Item(4,9).Value = 0.3910371682630009
$ws.Cells.Item(4,10).Value = 0.3910371682630009
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.229665
$ws.Cells.Item(4,14).Value = 6.688995
$ws.Cells.Item(4,15).Value = 0.09011208596621291
$ws.Cells.Item(4,16).Value = 0.09011208596621292
$ws.Cells.Item(4,17).Value = 362.096028545505
$ws.Cells.Item(4,18).Value = 3258.864256909545
$ws.Cells.Item(4,19).Value = 0.0352371749225
$ws.Cells.Item(4,20).Value = 0.03523717492250001

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 162.399297
$ws.Cells.Item(5,8).Value = 487.197891
$ws.Cells.Item(5,9).Value = 0.3910371682630009
$ws.Cells.Item(5,10).Value = 0.3910371682630009
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.074542333333333
$ws.Cells.Item(5,14).Value = 3.223627
$ws.Cells.Item(5,15).Value = 0.04342771273517247
$ws.Cells.Item(5,16).Value = 0.04342771273517248
$ws.Cells.Item(5,17).Value = 174.504919530073
$ws.Cells.Item(5,18).Value = 1570.544275770657
$ws.Cells.Item(5,19).Value = 0.01698184981210091
$ws.Cells.Item(5,20).Value = 0.01698184981210091

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 65.41736466666667
$ws.Cells.Item(6,8).Value = 196.252094
$ws.Cells.Item(6,9).Value = 0.1575168212364948
$ws.Cells.Item(6,10).Value = 0.1575168212364948
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 21.38398133333333
$ws.Cells.Item(6,14).Value = 64.151944
$ws.Cells.Item(6,15).Value = 0.864235283869651
$ws.Cells.Item(6,16).Value = 0.8642352838696511
$ws.Cells.Item(6,17).Value = 1398.88370490786
$ws.Cells.Item(6,18).Value = 12589.95334417074
$ws.Cells.Item(6,19).Value = 0.1361315947155672
$ws.Cells.Item(6,20).Value = 0.1361315947155672

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 65.41736466666667
$ws.Cells.Item(7,8).Value = 196.252094
$ws.Cells.Item(7,9).Value = 0.1575168212364948
$ws.Cells.Item(7,10).Value = 0.1575168212364948
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.05505166666666667
$ws.Cells.Item(7,14).Value = 0.165155
$ws.Cells.Item(7,15).Value = 0.002224917428963528
$ws.Cells.Item(7,16).Value = 0.002224917428963528
$ws.Cells.Item(7,17).Value = 3.601334953841111
$ws.Cells.Item(7,18).Value = 32.41201458457
$ws.Cells.Item(7,19).Value = 0.0003504619209240097
$ws.Cells.Item(7,20).Value = 0.0003504619209240097

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 65.41736466666667
$ws.Cells.Item(8,8).Value = 196.252094
$ws.Cells.Item(8,9).Value = 0.1575168212364948
$ws.Cells.Item(8,10).Value = 0.1575168212364948
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.229665
$ws.Cells.Item(8,14).Value = 6.688995
$ws.Cells.Item(8,15).Value = 0.09011208596621291
$ws.Cells.Item(8,16).Value = 0.09011208596621292
$ws.Cells.Item(8,17).Value = 145.8588083895034
$ws.Cells.Item(8,18).Value = 1312.72927550553
$ws.Cells.Item(8,19).Value = 0.01419416933638761
$ws.Cells.Item(8,20).Value = 0.01419416933638762

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 65.41736466666667
$ws.Cells.Item(9,8).Value = 196.252094
$ws.Cells.Item(9,9).Value = 0.1575168212364948
$ws.Cells.Item(9,10).Value = 0.1575168212364948
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.074542333333333
$ws.Cells.Item(9,14).Value = 3.223627
$ws.Cells.Item(9,15).Value = 0.04342771273517247
$ws.Cells.Item(9,16).Value = 0.04342771273517248
$ws.Cells.Item(9,17).Value = 70.29372766943756
$ws.Cells.Item(9,18).Value = 632.643549024938
$ws.Cells.Item(9,19).Value = 0.006840595263616012
$ws.Cells.Item(9,20).Value = 0.006840595263616013

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 126.3069433333333
$ws.Cells.Item(10,8).Value = 378.92083
$ws.Cells.Item(10,9).Value = 0.3041313008456065
$ws.Cells.Item(10,10).Value = 0.3041313008456065
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 21.38398133333333
$ws.Cells.Item(10,14).Value = 64.151944
$ws.Cells.Item(10,15).Value = 0.864235283869651
$ws.Cells.Item(10,16).Value = 0.8642352838696511
$ws.Cells.Item(10,17).Value = 2700.945318510391
$ws.Cells.Item(10,18).Value = 24308.50786659352
$ws.Cells.Item(10,19).Value = 0.2628410011199489
$ws.Cells.Item(10,20).Value = 0.262841001119949

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 126.3069433333333
$ws.Cells.Item(11,8).Value = 378.92083
$ws.Cells.Item(11,9).Value = 0.3041313008456065
$ws.Cells.Item(11,10).Value = 0.3041313008456065
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.05505166666666667
$ws.Cells.Item(11,14).Value = 0.165155
$ws.Cells.Item(11,15).Value = 0.002224917428963528
$ws.Cells.Item(11,16).Value = 0.002224917428963528
$ws.Cells.Item(11,17).Value = 6.953407742072223
$ws.Cells.Item(11,18).Value = 62.58066967865
$ws.Cells.Item(11,19).Value = 0.0006766670319447399
$ws.Cells.Item(11,20).Value = 0.0006766670319447402

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 126.3069433333333
$ws.Cells.Item(12,8).Value = 378.92083
$ws.Cells.Item(12,9).Value = 0.3041313008456065
$ws.Cells.Item(12,10).Value = 0.3041313008456065
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.229665
$ws.Cells.Item(12,14).Value = 6.688995
$ws.Cells.Item(12,15).Value = 0.09011208596621291
$ws.Cells.Item(12,16).Value = 0.09011208596621292
$ws.Cells.Item(12,17).Value = 281.6221708073167
$ws.Cells.Item(12,18).Value = 2534.59953726585
$ws.Cells.Item(12,19).Value = 0.02740590592681545
$ws.Cells.Item(12,20).Value = 0.02740590592681546

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 126.3069433333333
$ws.Cells.Item(13,8).Value = 378.92083
$ws.Cells.Item(13,9).Value = 0.3041313008456065
$ws.Cells.Item(13,10).Value = 0.3041313008456065
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.074542333333333
$ws.Cells.Item(13,14).Value = 3.223627
$ws.Cells.Item(13,15).Value = 0.04342771273517247
$ws.Cells.Item(13,16).Value = 0.04342771273517248
$ws.Cells.Item(13,17).Value = 135.7221576056011
$ws.Cells.Item(13,18).Value = 1221.49941845041
$ws.Cells.Item(13,19).Value = 0.01320772676689731
$ws.Cells.Item(13,20).Value = 0.01320772676689732

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Gnai2"
$ws.Cells.Item(14,3).Value = "Adra2b"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 61.180387
$ws.Cells.Item(14,8).Value = 183.541161
$ws.Cells.Item(14,9).Value = 0.1473147096548978
$ws.Cells.Item(14,10).Value = 0.1473147096548978
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 21.38398133333333
$ws.Cells.Item(14,14).Value = 64.151944
$ws.Cells.Item(14,15).Value = 0.864235283869651
$ws.Cells.Item(14,16).Value = 0.8642352838696511
$ws.Cells.Item(14,17).Value = 1308.280253574109
$ws.Cells.Item(14,18).Value = 11774.52228216698
$ws.Cells.Item(14,19).Value = 0.1273145699167758
$ws.Cells.Item(14,20).Value = 0.1273145699167758

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Gnai2"
$ws.Cells.Item(15,3).Value = "Adra2b"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 61.180387
$ws.Cells.Item(15,8).Value = 183.541161
$ws.Cells.Item(15,9).Value = 0.1473147096548978
$ws.Cells.Item(15,10).Value = 0.1473147096548978
$ws.Cells.Item(15,11).Value = 1
$ws.Cells.Item(15,12).Value = 0.3333333333333333
$ws.Cells.Item(15,13).Value = 0.05505166666666667
$ws.Cells.Item(15,14).Value = 0.165155
$ws.Cells.Item(15,15).Value = 0.002224917428963528
$ws.Cells.Item(15,16).Value = 0.002224917428963528
$ws.Cells.Item(15,17).Value = 3.368082271661666
$ws.Cells.Item(15,18).Value = 30.312740444955
$ws.Cells.Item(15,19).Value = 0.0003277630650538838
$ws.Cells.Item(15,20).Value = 0.0003277630650538839

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Gnai2"
$ws.Cells.Item(16,3).Value = "Adra2b"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 61.180387
$ws.Cells.Item(16,8).Value = 183.541161
$ws.Cells.Item(16,9).Value = 0.1473147096548978
$ws.Cells.Item(16,10).Value = 0.1473147096548978
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 2.229665
$ws.Cells.Item(16,14).Value = 6.688995
$ws.Cells.Item(16,15).Value = 0.09011208596621291
$ws.Cells.Item(16,16).Value = 0.09011208596621292
$ws.Cells.Item(16,17).Value = 136.411767580355
$ws.Cells.Item(16,18).Value = 1227.705908223195
$ws.Cells.Item(16,19).Value = 0.01327483578050984
$ws.Cells.Item(16,20).Value = 0.01327483578050985

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Gnai2"
$ws.Cells.Item(17,3).Value = "Adra2b"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 61.180387
$ws.Cells.Item(17,8).Value = 183.541161
$ws.Cells.Item(17,9).Value = 0.1473147096548978
$ws.Cells.Item(17,10).Value = 0.1473147096548978
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 1.074542333333333
$ws.Cells.Item(17,14).Value = 3.223627
$ws.Cells.Item(17,15).Value = 0.04342771273517247
$ws.Cells.Item(17,16).Value = 0.04342771273517248
$ws.Cells.Item(17,17).Value = 65.74091580121633
$ws.Cells.Item(17,18).Value = 591.6682422109469
$ws.Cells.Item(17,19).Value = 0.00639754089255824
$ws.Cells.Item(17,20).Value = 0.006397540892558242
